# chore: adapt column header formatting to respective input file names
#
# Rename the "_old"/"_new" header suffixes (columns A:J and L:U, row 1) to
# "_FV2304"/"_FV2310" respectively, turn the used range A1:U82 into a real
# Excel Table ("Table1") and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

# Columns A..J (1..10) carried the "_old" suffix
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304Headers[$i]
}

# Column K (11) is "diff" - untouched

# Columns L..U (12..21) carried the "_new" suffix
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2310Headers[$i]
}

# Turn A1:U82 into an Excel Table with a header row, named "Table1"
$tableRange = $ws.Range("A1:U82")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split/freeze pane after row 1)
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
